# Remove the last two slides (the "Rotation" / "No Rotation" demo slides,
# slide6.xml and slide7.xml) from the presentation.
$p = $ppt.ActivePresentation

# Delete slides from the end so indices of the remaining slides stay valid.
$p.Slides.Item($p.Slides.Count).Delete()
$p.Slides.Item($p.Slides.Count).Delete()
